$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): add P1=14, Q1=15, matching the formatting of O1 ---
$o1 = $ws.Range("O1")
$p1 = $ws.Range("P1")
$q1 = $ws.Range("Q1")
$o1.Copy($p1)
$o1.Copy($q1)
$p1.Value = 14
$q1.Value = 15

# --- Rows 2-25: swap I/K and M/O columns, then append P and Q columns with value 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # column K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # column M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # column O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # column P: new, value 2
    $ws.Cells.Item($r, 17).Value = 2  # column Q: new, value 2
}
